# HELLO / WORLD doc edit
#
#   Paragraph 1: "HELLO" -> "HLL"
#   Paragraph 2: "WORLD" -> "W" + bookmark(_GoBack) + "RL" + "D"
#
# The existing "_GoBack" bookmark currently sits at the very end of
# paragraph 2 (after the "D" run). After the edit it needs to sit right
# after the new "W" run, i.e. between "W" and "RL".

$d = $word.ActiveDocument

# 1) HELLO -> HLL  (plain text substitution, formatting is untouched)
$rHello = $d.Content
$rHello.Find.Execute("HELLO", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "HLL", 2)

# 2) Drop the current "_GoBack" bookmark; it will be re-created below in
#    its new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3) WORL -> W  (shrinks the first run of paragraph 2 down to just "W")
$rWorl = $d.Content
$rWorl.Find.Execute("WORL", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "W", 2)

# 4) Insert "RL" immediately after the "W" we just produced. Collapse the
#    found range to its end first so InsertAfter lands right after "W"
#    (and before the trailing "D" run).
$rWorl.Collapse(0)
$insertPos = $rWorl.Start
$rWorl.InsertAfter("RL")

# 5) Re-create the "_GoBack" bookmark, collapsed, at the point that sits
#    between "W" and "RL". Build a brand-new Range from the saved
#    character offset rather than reusing $rWorl, since InsertAfter grew
#    $rWorl to span the newly-inserted "RL" text.
$bookmarkRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
